# Creating noun declension table
# - Remove the stray "Α6 / ο, pl. -αε" row (old row 8); everything below
#   shifts up one row.
# - Append three new summary rows for the Α (masculine), Θ (feminine) and
#   Υ (neuter) declension groups, extending the Tabla2 table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the old row 8 ("Α6" / "ο, pl. -αε"); rows 9-30 shift up to 8-29.
$ws.Rows.Item(8).Delete()

# 2) Grow the table (ListObject) so it covers the two brand-new rows we are
#    about to add (the old blank trailing row, now row 29, plus two more).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B31"))

# 3) Fill in the three summary rows.
$ws.Range("A29").Value = "Α"
$ws.Range("B29").Value = "ο"
$ws.Range("A30").Value = "Θ"
$ws.Range("B30").Value = "α"
$ws.Range("A31").Value = "Υ"
$ws.Range("B31").Value = "το"

# Match the look of the rest of the table (14pt font, same as every other
# data row in the sheet).
$ws.Range("A29:B31").Font.Size = 14
$ws.Rows.Item(31).RowHeight = 18.75

# 4) Restore a sane view: no frozen/scrolled top-left cell, selection on A4.
$ws.Range("A4").Select() | Out-Null
